$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.280.00"
$ws.Range("E2").Value = "  -0.30%  "

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.630.95"
$ws.Range("E3").Value = "  -2.09%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.48"
$ws.Range("E5").Value = "  -3.10%  "

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.03"
$ws.Range("E6").Value = "  -1.53%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.02%  "

# Row 8: XRP
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.519"
$ws.Range("E8").Value = "  -1.16%  "

# Row 9: LidoStakedEther
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.631.26"
$ws.Range("E9").Value = "  -2.04%  "

# Row 10: Dogecoin
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.171"
$ws.Range("E10").Value = "  +0.06%  "

# Row 11: TRON
$ws.Range("E11").Value = "  +1.28%  "

# Row 12: Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.357"
$ws.Range("E12").Value = "  +0.38%  "

# Row 13: Toncoin
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.92"
$ws.Range("E13").Value = "  -2.37%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.110.94"
$ws.Range("E14").Value = "  -1.18%  "

# Row 15: ShibaInu
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000185"
$ws.Range("E15").Value = "  -0.84%  "

# Row 16: WrappedBTC
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "72.184.63"
$ws.Range("E16").Value = "  -0.21%  "

# Row 17: Avalanche
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.72"
$ws.Range("E17").Value = "  -2.69%  "

# Row 18: WrappedEther
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.624.03"
$ws.Range("E18").Value = "  -2.24%  "

# Row 19: Chainlink
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.03"
$ws.Range("E19").Value = "  -0.64%  "

# Row 20: Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.86"
$ws.Range("E20").Value = "  -1.14%  "

# Row 21: BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.63"
$ws.Range("E21").Value = "  +0.41%  "

# Row 22: Polkadot
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.10"
$ws.Range("E22").Value = "  -2.03%  "

# Row 23: SuiNetwork
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.03"
$ws.Range("E23").Value = "  -1.28%  "

# Row 24: Litecoin
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.25"
$ws.Range("E24").Value = "  -1.50%  "

# Row 25: Dai
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.21%  "

# Row 26: NEARProtocol
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.22"
$ws.Range("E26").Value = "  -3.23%  "

# Row 27: Aptos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.45"
$ws.Range("E27").Value = "  -4.24%  "

# Row 28: WrappedeETH
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.768.06"
$ws.Range("E28").Value = "  -2.10%  "

# Row 29: Binance-PegBSC-USD
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.07%  "

# Row 30: PEPE
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0944"
$ws.Range("E30").Value = "  -0.38%  "

# Row 31: InternetComputer(DFINITY)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.91"
$ws.Range("E31").Value = "  -2.56%  "

# Row 32: Bittensor
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "492.35"
$ws.Range("E32").Value = "  -4.44%  "

# Row 33: Fetch.AI
$ws.Range("E33").Value = "  -2.77%  "

# Row 34: PancakeSwap
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.78"
$ws.Range("E34").Value = "  -2.18%  "

# Row 36: Monero
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.74"
$ws.Range("E36").Value = "  -2.17%  "

# Row 37: Kaspa
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.114"
$ws.Range("E37").Value = "  +4.75%  "

# Row 38: EthereumClassic
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.07"
$ws.Range("E38").Value = "  -2.90%  "

# Row 39: WhiteBITCoin
$ws.Range("E39").Value = "  -1.45%  "

# Row 40: ImmutableX
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.35"
$ws.Range("E40").Value = "  -2.40%  "

# Row 41: USDe
$ws.Range("E41").Value = "  -0.12%  "

# Row 42: Stacks
$ws.Range("E42").Value = "  -6.02%  "

# Row 43: dogwifhat
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.55"
$ws.Range("E43").Value = "  -1.48%  "

# Row 44: RenderToken
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.86"
$ws.Range("E44").Value = "  -3.67%  "

# Row 45: PolygonEcosystemToken
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.324"
$ws.Range("E45").Value = "  -2.95%  "

# Row 46: OKB
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.02"
$ws.Range("E46").Value = "  -0.68%  "

# Row 47: Aave
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.55"
$ws.Range("E47").Value = "  -2.25%  "

# Row 48: Filecoin
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.63"
$ws.Range("E48").Value = "  -3.20%  "

# Row 49: ARBITRUM
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.543"
$ws.Range("E49").Value = "  -2.06%  "

# Row 50: Optimism
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.66"
$ws.Range("E50").Value = "  -4.24%  "

# Row 51: Mantle
$ws.Range("E51").Value = "  -0.04%  "
